$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.535.96"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "1.915.59"
$ws.Range("E3").Value = "  +2.72%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5159"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.83%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3983"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09867"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.152"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.529"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.53%  "

$ws.Range("D14").Value = "1.911.80"
$ws.Range("E14").Value = "  +3.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.472"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001139"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.91%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.315"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.12%  "

$ws.Range("D23").Value = "28.595.87"
$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.680"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.97%  "

$ws.Range("D27").Value = "2.137.44"
$ws.Range("E27").Value = "  +3.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.63%  "

$ws.Range("E31").Value = "  +5.89%  "

$ws.Range("E32").Value = "  +1.81%  "

$ws.Range("E33").Value = "  +1.60%  "

$ws.Range("E34").Value = "  +1.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.838"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06766"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02443"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.273"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2225"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6487"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.096"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.189"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6111"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.766"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.286"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.064"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.206"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "
